# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# for rows 2-51 on Sheet1, matching the Fri Jul 28 14:49:19 UTC 2023
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new text value. A leading
# apostrophe forces Excel to store the value as text (matching the
# original inline-string cell type) instead of auto-converting
# numeric-looking strings (e.g. "242.14") into a Number cell; the
# Style reset afterwards keeps the cell on the default/unstyled xf
# so no stray number formatting is left behind.
$updates = @(
    @{ Cell = "D2"; Value = "29.423.38" }
    @{ Cell = "E2"; Value = "  -0.04%  " }
    @{ Cell = "D3"; Value = "1.875.60" }
    @{ Cell = "E3"; Value = "  +0.02%  " }
    @{ Cell = "E4"; Value = "  -0.34%  " }
    @{ Cell = "D5"; Value = "0.7143" }
    @{ Cell = "E5"; Value = "  +0.06%  " }
    @{ Cell = "D6"; Value = "242.14" }
    @{ Cell = "E6"; Value = "  +0.22%  " }
    @{ Cell = "D7"; Value = "0.9999" }
    @{ Cell = "E7"; Value = "  -0.24%  " }
    @{ Cell = "D8"; Value = "0.3116" }
    @{ Cell = "E8"; Value = "  +0.94%  " }
    @{ Cell = "D9"; Value = "0.07736" }
    @{ Cell = "E9"; Value = "  -2.04%  " }
    @{ Cell = "D10"; Value = "25.02" }
    @{ Cell = "E10"; Value = "  -1.68%  " }
    @{ Cell = "D11"; Value = "0.08379" }
    @{ Cell = "E11"; Value = "  +1.53%  " }
    @{ Cell = "D12"; Value = "1.906.51" }
    @{ Cell = "E12"; Value = "  +1.94%  " }
    @{ Cell = "D13"; Value = "5.247" }
    @{ Cell = "E13"; Value = "  -0.15%  " }
    @{ Cell = "D14"; Value = "0.7186" }
    @{ Cell = "E14"; Value = "  -0.80%  " }
    @{ Cell = "D15"; Value = "91.55" }
    @{ Cell = "E15"; Value = "  +0.46%  " }
    @{ Cell = "D16"; Value = "29.429.09" }
    @{ Cell = "E16"; Value = "  +0.05%  " }
    @{ Cell = "D17"; Value = "0.000008197" }
    @{ Cell = "E17"; Value = "  +4.74%  " }
    @{ Cell = "D18"; Value = "5.978" }
    @{ Cell = "E18"; Value = "  +1.95%  " }
    @{ Cell = "E19"; Value = "  -0.06%  " }
    @{ Cell = "D20"; Value = "2.128.00" }
    @{ Cell = "E20"; Value = "  +1.13%  " }
    @{ Cell = "D21"; Value = "13.22" }
    @{ Cell = "E21"; Value = "  -0.05%  " }
    @{ Cell = "D22"; Value = "0.9993" }
    @{ Cell = "E22"; Value = "  -0.26%  " }
    @{ Cell = "D23"; Value = "7.941" }
    @{ Cell = "E23"; Value = "  -1.31%  " }
    @{ Cell = "D24"; Value = "0.9993" }
    @{ Cell = "E24"; Value = "  -0.40%  " }
    @{ Cell = "D25"; Value = "0.1631" }
    @{ Cell = "E25"; Value = "  +1.82%  " }
    @{ Cell = "D26"; Value = "163.81" }
    @{ Cell = "E26"; Value = "  +0.71%  " }
    @{ Cell = "D27"; Value = "9.035" }
    @{ Cell = "E27"; Value = "  +0.33%  " }
    @{ Cell = "D28"; Value = "18.57" }
    @{ Cell = "E28"; Value = "  +1.50%  " }
    @{ Cell = "D29"; Value = "1.509" }
    @{ Cell = "E29"; Value = "  +0.77%  " }
    @{ Cell = "E30"; Value = "  +0.44%  " }
    @{ Cell = "D31"; Value = "1.299" }
    @{ Cell = "E31"; Value = "  -4.25%  " }
    @{ Cell = "D32"; Value = "4.328" }
    @{ Cell = "E32"; Value = "  +5.60%  " }
    @{ Cell = "D33"; Value = "0.05246" }
    @{ Cell = "E33"; Value = "  +1.02%  " }
    @{ Cell = "D34"; Value = "1.932" }
    @{ Cell = "E34"; Value = "  -0.04%  " }
    @{ Cell = "D35"; Value = "0.7682" }
    @{ Cell = "E35"; Value = "  +6.44%  " }
    @{ Cell = "E36"; Value = "  -1.48%  " }
    @{ Cell = "D37"; Value = "2.685" }
    @{ Cell = "E37"; Value = "  +0.37%  " }
    @{ Cell = "E38"; Value = "  +0.27%  " }
    @{ Cell = "D39"; Value = "2.723" }
    @{ Cell = "E39"; Value = "  +1.25%  " }
    @{ Cell = "D40"; Value = "1.167.93" }
    @{ Cell = "E40"; Value = "  -1.03%  " }
    @{ Cell = "D41"; Value = "6.428" }
    @{ Cell = "E41"; Value = "  +4.90%  " }
    @{ Cell = "D42"; Value = "73.66" }
    @{ Cell = "E42"; Value = "  +1.60%  " }
    @{ Cell = "D43"; Value = "0.8915" }
    @{ Cell = "E43"; Value = "  -1.85%  " }
    @{ Cell = "D44"; Value = "104.05" }
    @{ Cell = "E44"; Value = "  +1.87%  " }
    @{ Cell = "D45"; Value = "0.9991" }
    @{ Cell = "D46"; Value = "2.024.38" }
    @{ Cell = "E46"; Value = "  +0.60%  " }
    @{ Cell = "D47"; Value = "1.807" }
    @{ Cell = "E47"; Value = "  +0.86%  " }
    @{ Cell = "D48"; Value = "0.5201" }
    @{ Cell = "E48"; Value = "  -1.82%  " }
    @{ Cell = "D49"; Value = "9.420" }
    @{ Cell = "E49"; Value = "  +1.41%  " }
    @{ Cell = "E50"; Value = "  +0.67%  " }
    @{ Cell = "E51"; Value = "  +0.49%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
    $ws.Range($u.Cell).Style = "Normal"
}
